$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 3, shifting existing rows 3-10 down to 4-11
$ws.Rows("3:3").Insert()

# Populate the newly inserted row 3 with the new review data
$ws.Range("A3").Value = 5
$ws.Range("B3").Value = ""
$ws.Range("C3").Value = 45902.89908415509
$ws.Range("D3").Value = "MzM3NmJjMjQtMTRmMi00NjgxLWI4ZmQtZjJmNjQwMjkxNzhkOjU3MDE2"
